$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/6/2023  Through  2/12/2023"

# --- Style/type conversions: text<->numeric cells (copy format+value from a donor cell sharing the target style,
#     then overwrite the value where it must differ from the donor) ---
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("C14").Copy($ws.Range("C22"))
$ws.Range("C36").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("K36").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100
$ws.Range("K36").Copy($ws.Range("L22"))
$ws.Range("L22").Value = 0
$ws.Range("C14").Copy($ws.Range("C26"))
$ws.Range("C14").Copy($ws.Range("D26"))
$ws.Range("E14").Copy($ws.Range("E26"))
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
$ws.Range("C36").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 1
$ws.Range("K36").Copy($ws.Range("M28"))
$ws.Range("M28").Value = 100
$ws.Range("C36").Copy($ws.Range("C29"))
$ws.Range("C29").Value = 1
$ws.Range("K36").Copy($ws.Range("M29"))
$ws.Range("M29").Value = 100

$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("N15").Value = -44.444444444444
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -25
$ws.Range("F16").Value = 19
$ws.Range("G16").Value = 30
$ws.Range("H16").Value = -36.666666666666
$ws.Range("I16").Value = 35
$ws.Range("J16").Value = 48
$ws.Range("K16").Value = -27.083333333333
$ws.Range("L16").Value = 20.689655172413
$ws.Range("M16").Value = -37.5
$ws.Range("N16").Value = -82.673267326732
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 83.333333333333
$ws.Range("F17").Value = 46
$ws.Range("G17").Value = 33
$ws.Range("H17").Value = 39.393939393939
$ws.Range("I17").Value = 65
$ws.Range("J17").Value = 57
$ws.Range("K17").Value = 14.035087719298
$ws.Range("L17").Value = 47.727272727272
$ws.Range("M17").Value = 116.666666666667
$ws.Range("N17").Value = -26.966292134831
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 66.666666666666
$ws.Range("F18").Value = 14
$ws.Range("H18").Value = 40
$ws.Range("I18").Value = 25
$ws.Range("J18").Value = 18
$ws.Range("K18").Value = 38.888888888888
$ws.Range("L18").Value = 92.307692307692
$ws.Range("M18").Value = -10.714285714285
$ws.Range("N18").Value = -81.751824817518
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -35.714285714285
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = -21.276595744680
$ws.Range("I19").Value = 60
$ws.Range("J19").Value = 77
$ws.Range("K19").Value = -22.077922077922
$ws.Range("L19").Value = 46.341463414634
$ws.Range("M19").Value = 22.448979591836
$ws.Range("N19").Value = -53.488372093023
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -20
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 22
$ws.Range("H20").Value = -18.181818181818
$ws.Range("I20").Value = 29
$ws.Range("J20").Value = 40
$ws.Range("K20").Value = -27.5
$ws.Range("L20").Value = 81.25
$ws.Range("M20").Value = 52.631578947368
$ws.Range("N20").Value = -81.168831168831
$ws.Range("D21").Value = 36
$ws.Range("E21").Value = -2.777777777777
$ws.Range("F21").Value = 136
$ws.Range("G21").Value = 144
$ws.Range("H21").Value = -5.555555555555
$ws.Range("I21").Value = 219
$ws.Range("J21").Value = 245
$ws.Range("K21").Value = -10.612244897959
$ws.Range("L21").Value = 51.034482758620
$ws.Range("M21").Value = 18.378378378378
$ws.Range("N21").Value = -69.667590027700
$ws.Range("J22").Value = 4
$ws.Range("K22").Value = -75
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -50
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 33.333333333333
$ws.Range("I23").Value = 6
$ws.Range("J23").Value = 6
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 100
$ws.Range("M23").Value = -25
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = -9.090909090909
$ws.Range("F24").Value = 128
$ws.Range("G24").Value = 101
$ws.Range("H24").Value = 26.732673267326
$ws.Range("I24").Value = 189
$ws.Range("J24").Value = 164
$ws.Range("K24").Value = 15.243902439024
$ws.Range("L24").Value = 56.198347107438
$ws.Range("M24").Value = 58.823529411764
$ws.Range("C25").Value = 17
$ws.Range("D25").Value = 21
$ws.Range("E25").Value = -19.047619047619
$ws.Range("F25").Value = 60
$ws.Range("G25").Value = 59
$ws.Range("H25").Value = 1.694915254237
$ws.Range("I25").Value = 104
$ws.Range("J25").Value = 84
$ws.Range("K25").Value = 23.809523809523
$ws.Range("L25").Value = 136.363636363636
$ws.Range("M25").Value = 2.970297029702
$ws.Range("F26").Value = 4
$ws.Range("H26").Value = 33.333333333333
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 0
$ws.Range("L27").Value = -14.285714285714
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 2
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 100
$ws.Range("N28").Value = -81.818181818181
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 2
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 100
$ws.Range("N29").Value = -80
